# Regenerated BOM outputs after a BOM change: several capacitor and resistor
# designator groups were re-bucketed (C200 moved out of its own row into the
# 1nF group with C301-C303; C300/C304 now form their own 10nF row; the 270k
# and 10k resistor groups swapped several designators; the 2k and 100R groups
# lost a couple of designators to the 270k/10k groups). Quantities are
# updated to match the new designator counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => [Quantity, Designator, Type, Comment, Description, Manufacturer,
#         ManufacturerPartNumber, LibRef, Supplier1, SupplierPartNumber1]
$rows = @{
    2  = @(2,  'C300, C304', 'Capacitor', 'C 10nF 0402 16V', 'CAP CER 10000PF 16V Y5V 0402', 'Samsung Electro-Mechanics', 'CL05F103ZO5NNNC', 'CAP_10nF_16V_0402', 'Digi-Key', '1276-1738-1-ND')
    3  = @(2,  'C307, C308', 'Capacitor', 'C 5.6pF 0402 50V', 'CAP CER 5.6PF 50V C0G/NP0 0402', 'Samsung Electro-Mechanics', 'CL05C5R6DB5NNNC', 'CAP_5.6PF_50V_0402', 'Digi-Key', '1276-1712-1-ND')
    4  = @(4,  'C200, C301, C302, C303', 'Capacitor', 'C 1nF 0402 50V', 'CAP CER 1000PF 50V C0G/NP0 0402', 'Murata Electronics', 'GRM1555C1H102JA01J', 'CAP_1nF_50V_0402', 'Digi-Key', '490-6190-1-ND')
    5  = @(4,  'C201, C202, C203, C208', 'Capacitor', 'C 10uF 0603 6.3V', 'CAP CER 10UF 6.3V X5R 0603', 'Taiyo Yuden', 'JMK107ABJ106MA-T', 'CAP_10uF_6V3_0603', 'Digi-Key', '587-5869-1-ND')
    6  = @(8,  'C204, C205, C206, C207, C209, C305, C306, C400', 'Capacitor', 'C 100nF 0402 10V', 'CAP CER 0.1UF 10V X5R 0402', 'Samsung Electro-Mechanics', 'CL05A104MP5NNNC', 'CAP_100nF_10V_0402', 'Digi-Key', '1276-1443-1-ND')
    37 = @(13, 'R200, R201, R216, R301, R302, R303, R304, R305, R307, R310, R315, R319, R321', 'Resistor', 'R 270k 0402', 'RES 270K OHM 1% 1/16W 0402', 'YAGEO', 'RC0402FR-07270KL', 'R_270k_0402', 'Digi-Key', '311-270KLRCT-ND')
    38 = @(13, 'R203, R205, R208, R209, R210, R214, R215, R300, R306, R312, R313, R314, R320', 'Resistor', 'R 10k 0402', 'RES 10K OHM 1% 1/16W 0402', 'YAGEO', 'RC0402FR-0710KL', 'R_10k_0402', 'Digi-Key', '311-10.0KLRCT-ND')
    39 = @(2,  'R211, R212', 'Resistor', 'R 5k1 0402', 'RES 5.1K OHM 1% 1/16W 0402', 'Stackpole Electronics Inc', 'RMCF0402FT5K10', 'R_5k1_0420', 'Digi-Key', 'RMCF0402FT5K10CT-ND')
    40 = @(3,  'R218, R219, R402', 'Resistor', 'R 0R0 0402', 'RES 0 OHM JUMPER 1/16W 0402', 'Stackpole Electronics Inc', 'RMCF0402ZT0R00', 'R_0R0_0402', 'Digi-Key', 'RMCF0402ZT0R00CT-ND')
    41 = @(6,  'R204, R206, R207, R308, R309, R311', 'Resistor', 'R 2k 0402', 'RES 2K OHM 1% 1/16W 0402', 'Stackpole Electronics Inc', 'RMCF0402FT2K00', 'R_2k_0402', 'Digi-Key', 'RMCF0402FT2K00CT-ND')
    42 = @(7,  'R202, R316, R317, R318, R322, R323, R324', 'Resistor', 'R 100R 0402', 'RES 100 OHM 1% 1/16W 0402', 'Stackpole Electronics Inc', 'RMCF0402FT100R', 'R_100R_0402', 'Digi-Key', 'RMCF0402FT100RCT-ND')
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $values[$col - 1]
    }
}
